$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Delete the rows (by account number, matching column A) that were removed
# in the target revision. Working from the bottom up keeps earlier row
# indices valid as we go.
$rowsToDelete = @(9, 8, 6, 4, 3, 2)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
